$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- new content (previously row 4's data)
$ws.Range("A3").Value = 112436438
$ws.Range("B3").Value = 56350
$ws.Range("E3").Value = 102110
$ws.Range("F3").Value = "Fjällvråk"
$ws.Range("G3").Value = "Buteo lagopus"
$ws.Range("H3").Value = "(Pontoppidan, 1763)"
$ws.Range("I3").Value = "'1"
$ws.Range("M3").Value = "sträckande S"
$ws.Range("AC3").ClearContents()

# Row 4 <- new content (previously row 5's data)
$ws.Range("A4").Value = 112436408
$ws.Range("B4").Value = 57001
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 102120
$ws.Range("F4").Value = "Rödstrupig piplärka"
$ws.Range("G4").Value = "Anthus cervinus"
$ws.Range("H4").Value = "(Pallas, 1811)"
$ws.Range("I4").Value = "'2"
$ws.Range("M4").Value = "lockläte, övriga läten"

# Row 5 <- new content (previously row 3's data, with B changed)
$ws.Range("A5").Value = 112436368
$ws.Range("B5").Value = 56332
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100034
$ws.Range("F5").Value = "Blå kärrhök"
$ws.Range("G5").Value = "Circus cyaneus"
$ws.Range("H5").Value = "(Linnaeus, 1766)"
$ws.Range("M5").Value = "födosökande"
$ws.Range("AC5").Value = "Hans * hona."
